# Edit script for CubeA-HW20.xlsx
# Commit message: "Updated notebook, reran simulation"
#
# The underlying HW (homework-average) simulation was rerun with two
# additional trials ("Holden" and "Rizzie Spiral") inserted right after
# "Spiral5". Those become rows 4 and 5 on the sheet; everything that used
# to live in rows 4-29 shifts down two rows (now rows 6-31), and the two
# final trials (formerly rows 28-29, "Michael-CCHex"/"Michael-SNHex") are
# duplicated onto brand-new rows 30 and 31. One method name was also
# corrected: "Thomas Hex" -> "Matthies Hex" (row 9). The sheet's used
# range grows from A1:W29 to A1:W31.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Rename "Thomas Hex" -> "Matthies Hex" -------------------------------------------
$ws.Cells.Item(9, 2).Value2 = "Matthies Hex"

# --- 2. Row labels (column A = trial index, column B = trial name) for rows that change ---
$rowLabels = @{
    4  = @(2,  "Holden")
    5  = @(3,  "Rizzie Spiral")
    30 = @(28, "Michael-CCHex")
    31 = @(29, "Michael-SNHex")
}
foreach ($r in $rowLabels.Keys) {
    $pair = $rowLabels[$r]
    $ws.Cells.Item($r, 1).Value2 = $pair[0]
    $ws.Cells.Item($r, 2).Value2 = $pair[1]
}

# --- 3. Numeric data block, columns C:W (3..23), rows 4-31 -------------------------------
$rows = @(
    @(1.037696498706205, 1.37061559928647, 1.37061559928647, 0.8967181020725841, 0.8419190577247671, 0.8419190577247671, 0.8419190577247671, 0.880961677183441, 1.054467251998644, 0.9244589414005867, 1.213892286069556, 0.8419190577247671, 0.8967181020725841, 1.133666850679527, 0.975592677035614, 1.036417586361274, 1.107266984452566, 1.036417586361274, 1.040930002770616, 1.001127813761447, 1.027591176805282),
    @(1.21237980931583, 1.765550952476064, 1.765550952476064, 1.021184500117918, 0.346056142105057, 0.346056142105057, 0.346056142105057, 0.7894356102796086, 1.095223269890927, 0.7103156254172743, 1.462304289233247, 0.346056142105057, 1.021184500117918, 1.393367726296991, 1.058203885004423, 1.04426386489968, 1.293986240828303, 1.04426386489968, 1.057003716147491, 0.9148142013390046, 1.050306274854491),
    @(1.042115245659225, 1.36870651348404, 1.36870651348404, 0.9007147947639962, 0.8322912771900883, 0.8322912771900883, 0.8322912771900883, 0.8803710178967409, 1.055906010167462, 0.9216811425341112, 1.215365410615678, 0.8322912771900883, 0.9007147947639962, 1.134710654124018, 0.9783104024657293, 1.033904195146041, 1.108442439471833, 1.033904195146042, 1.039404648901397, 0.9979819745591352, 1.027143926538918),
    @(1.047374604884727, 1.461995371485593, 1.461995371485593, 0.8798063211023061, 0.8010647349553307, 0.8010647349553307, 0.8010647349553307, 0.855893715367433, 1.063723666822767, 0.9028590499567732, 1.263500918170029, 0.8010647349553307, 0.8798063211023061, 1.17090084629395, 0.9717649939625368, 1.04762214251441, 1.135175119803556, 1.04762214251441, 1.051647523591499, 1.001530965864265, 1.03452729784312),
    @(1.097130482808417, 2.656109404141916, 2.656109404141916, 0.5568150880210136, 0.4483900230707777, 0.4483900230707777, 0.4483900230707777, 0.524904094014556, 1.178596350589907, 0.6868586105271516, 1.887336381531393, 0.4483900230707777, 0.5568150880210136, 1.606462246081465, 0.8677057193054601, 1.220438171744569, 1.463840280917612, 1.220438171744569, 1.209977716455903, 1.057660177778878, 1.129517554338141),
    @(1.010184408263242, 1.042846022819864, 1.042846022819864, 0.9980358260533719, 0.9676516073423105, 0.9676516073423105, 0.9676516073423105, 0.9876801666157652, 1.005605799979156, 0.9855160528997359, 1.025677135731166, 0.9676516073423105, 0.9980358260533719, 1.020440924436618, 1.001820813016264, 1.002844485405182, 1.015495882950797, 1.002844485405182, 1.003534814048676, 0.9963581727074026, 1.002899627463077),
    @(1.000438000469342, 1.003839306369855, 1.003839306369855, 0.9994168266002679, 0.998128435043112, 0.998128435043112, 0.998128435043112, 0.998989343418796, 1.00038119888568, 0.9990310435735678, 1.002098984670025, 0.998128435043112, 0.9994168266002679, 1.001628066485061, 0.9998990127429738, 1.000461522671078, 1.001212443951934, 1.000461522671078, 1.000441441724728, 0.9999788403884052, 1.00029039237883),
    @(1.017024146827711, 1.072948661355545, 1.072948661355545, 0.997907526360858, 0.9455123443089511, 0.9455123443089511, 0.9455123443089511, 0.9798891519422757, 1.008630370789052, 0.975086319551327, 1.042902422397654, 0.9455123443089511, 0.997907526360858, 1.035428093858201, 1.003268948574955, 1.005456177341785, 1.026495519501818, 1.005456177341785, 1.006249725703601, 0.9941022494246713, 1.004987617941671),
    @(1.094784945641312, 2.69322854883295, 2.69322854883295, 0.5519893919522594, 0.4474462936349127, 0.4474462936349127, 0.4474462936349127, 0.5199906024143297, 1.175764032522311, 0.6794064581312521, 1.901431690981314, 0.4474462936349127, 0.5519893919522594, 1.622608970392605, 0.8638767122372852, 1.230888078140041, 1.473660657769173, 1.230888078140041, 1.217107066735608, 1.063174912115469, 1.13300524551383),
    @(1.000377852829018, 1.267020059974008, 1.267020059974008, 0.8898629817129452, 0.9577701891361196, 0.9577701891361196, 0.9577701891361196, 0.9129102717908916, 1.034438843493537, 0.9670852576667915, 1.144942954176899, 0.9577701891361196, 0.8898629817129452, 1.078441520843477, 0.962150912603241, 1.038217743607691, 1.06377396172683, 1.038217743607691, 1.037273018579153, 1.021372452690546, 1.021801051347526),
    @(1.031928635410527, 1.506451411185266, 1.506451411185266, 0.8672452017999999, 0.8320350566684199, 0.8320350566684199, 0.8320350566684199, 0.8547848749789464, 1.052912969105268, 0.9008243625157891, 1.272027712210529, 0.8320350566684199, 0.8672452017999999, 1.186848306492633, 0.9600790854526339, 1.068577223217895, 1.142203194030178, 1.068577223217895, 1.064661159689738, 1.018135939085475, 1.039776277984343),
    @(0.948415354163217, 0.7612809037377597, 0.7612809037377597, 1.084916004420995, 1.168697173486799, 1.168697173486799, 1.168697173486799, 1.104619725467156, 0.9312628468662161, 1.05146463038284, 0.8318362937857001, 1.168697173486799, 1.084916004420995, 0.9230984540793772, 1.008089425643605, 1.004964693881851, 0.9258199183416568, 1.004964693881851, 0.9865392321279425, 1.022970820399714, 0.9853116165388354),
    @(1.156870699999999, 3.939353100000003, 3.939353100000003, 0.2188589699999998, 0.05435712199999995, 0.05435712199999995, 0.05435712199999995, 0.1699082099999997, 1.3021676, 0.4490890699999996, 2.559413900000002, 0.05435712199999995, 0.2188589699999998, 2.079106035000001, 0.7605132850000001, 1.404189730666668, 1.820126556666668, 1.404189730666668, 1.378684198000001, 1.113818782800001, 1.231252334000001),
    @(1.1568707, 3.9393531, 3.9393531, 0.21885897, 0.054357122, 0.054357122, 0.054357122, 0.16990821, 1.3021676, 0.44908907, 2.5594139, 0.054357122, 0.21885897, 2.079106035, 0.760513285, 1.404189730666667, 1.820126556666667, 1.404189730666667, 1.378684198, 1.1138187828, 1.231252334),
    @(1.1568707, 3.9393531, 3.9393531, 0.21885897, 0.054357122, 0.054357122, 0.054357122, 0.16990821, 1.3021676, 0.44908907, 2.5594139, 0.054357122, 0.21885897, 2.079106035, 0.760513285, 1.404189730666667, 1.820126556666667, 1.404189730666667, 1.378684198, 1.1138187828, 1.231252334),
    @(1.240911, 0.82576034, 0.82576034, 1.3656747, 0.47174105, 0.47174105, 0.47174105, 1.0555159, 1.0027216, 0.8263528499999999, 0.98061785, 0.47174105, 1.3656747, 1.09571752, 1.18419815, 0.8877253633333333, 1.06471888, 0.8877253633333334, 0.9164744225000001, 0.8275277480000002, 0.9711619112500001),
    @(1.209545138356164, 1.886962695068492, 1.886962695068492, 0.9753982209589043, 0.3294725000821919, 0.3294725000821919, 0.3294725000821919, 0.7541493502739726, 1.107435649863014, 0.6958725647945209, 1.52421080630137, 0.3294725000821919, 0.9753982209589043, 1.431180458013698, 1.041416935410959, 1.06394447203653, 1.32326552196347, 1.063944472036529, 1.074817266493151, 0.9257483132109587, 1.060380865712329),
    @(1.207465884210527, 1.967978105789474, 1.967978105789474, 0.9455348294736843, 0.3186126825263157, 0.3186126825263157, 0.3186126825263157, 0.7310871999999999, 1.115122437368421, 0.6861311110526316, 1.565074086315789, 0.3186126825263157, 0.9455348294736843, 1.456756467631579, 1.030328633421053, 1.077375205929825, 1.34287845754386, 1.077375205929824, 1.086812013789474, 0.933172147536842, 1.067125792092105),
    @(1.207465884210527, 1.967978105789474, 1.967978105789474, 0.9455348294736843, 0.3186126825263157, 0.3186126825263157, 0.3186126825263157, 0.7310871999999999, 1.115122437368421, 0.6861311110526316, 1.565074086315789, 0.3186126825263157, 0.9455348294736843, 1.456756467631579, 1.030328633421053, 1.077375205929825, 1.34287845754386, 1.077375205929824, 1.086812013789474, 0.933172147536842, 1.067125792092105),
    @(1.097375869095959, 0.6474123980439753, 0.6474123980439753, 1.383247727733637, 0.7695134981280533, 0.7695134981280533, 0.7695134981280533, 1.176299178148569, 0.9210545075827427, 0.9294315679662116, 0.8036241904562451, 0.7695134981280533, 1.383247727733637, 1.015330062888806, 1.15215111765819, 0.9333912079685552, 0.983904877786785, 0.9333912079685552, 0.9303070328721021, 0.8981483259232924, 0.9659948671444241),
    @(0.9361305276650942, 1.013927747195857, 1.013927747195857, 0.8574694106388464, 1.165039530280905, 1.165039530280905, 1.165039530280905, 0.9650377884299971, 1.020185115454781, 1.069688241736297, 1.00572304026691, 1.165039530280905, 0.8574694106388464, 0.9356985789173518, 0.9388272630468135, 1.012145562705203, 0.9638607577631614, 1.012145562705203, 1.014155450892597, 1.044332266770259, 1.004150175208586),
    @(0.8745368999227625, 0.6132311431620122, 0.6132311431620122, 0.9489468280812738, 1.391676314260136, 1.391676314260136, 1.391676314260136, 1.097446714224897, 0.9525803801333577, 1.16504839110238, 0.7638635791204468, 1.391676314260136, 0.9489468280812738, 0.781088985621643, 0.9507636041073158, 0.9846180951678072, 0.8382527837922146, 0.9846180951678072, 0.9766086664091949, 1.059622195979383, 0.9759162812509081),
    @(0.9921446998649259, 0.9208543484117759, 0.9208543484117759, 1.114500948424145, 1.044247981991259, 1.044247981991259, 1.044247981991259, 1.074483415934369, 0.9330467897408178, 0.9797230111250645, 0.9170928555211185, 1.044247981991259, 1.114500948424145, 1.017677648417961, 1.023773869082481, 1.026534426275727, 0.9894673621922463, 1.026534426275727, 1.003162517142, 1.011379610111852, 0.9970117563766845),
    @(1.097375869095959, 0.6474123980439753, 0.6474123980439753, 1.383247727733637, 0.7695134981280533, 0.7695134981280533, 0.7695134981280533, 1.176299178148569, 0.9210545075827427, 0.9294315679662115, 0.8036241904562451, 0.7695134981280533, 1.383247727733637, 1.015330062888806, 1.15215111765819, 0.9333912079685552, 0.983904877786785, 0.9333912079685552, 0.9303070328721021, 0.8981483259232924, 0.9659948671444242),
    @(0.9361305276650941, 1.013927747195857, 1.013927747195857, 0.8574694106388463, 1.165039530280905, 1.165039530280905, 1.165039530280905, 0.9650377884299971, 1.020185115454781, 1.069688241736297, 1.00572304026691, 1.165039530280905, 0.8574694106388463, 0.9356985789173516, 0.9388272630468135, 1.012145562705203, 0.9638607577631614, 1.012145562705203, 1.014155450892597, 1.044332266770259, 1.004150175208586),
    @(1.040789098745243, 1.097862878032881, 1.097862878032881, 0.9594025300772111, 0.8953943628454462, 0.8953943628454462, 0.8953943628454462, 0.9429354979296316, 1.041540739109217, 0.9721315662460449, 1.083685311755372, 0.8953943628454462, 0.9594025300772111, 1.028632704055046, 1.000471634593214, 0.9842199236518461, 1.032935382406436, 0.9842199236518461, 0.9985501275161888, 0.9779189745820404, 1.004217748092631),
    @(1.062864580180446, 1.003375168059962, 1.003375168059962, 1.096089382620813, 0.8541696514648782, 0.8541696514648782, 0.8541696514648782, 1.008845290132268, 0.9963567386157495, 0.9412078236537136, 1.014056563744896, 0.8541696514648782, 1.096089382620813, 1.049732275340387, 1.046223060618281, 0.9845447340485508, 1.031940429765508, 0.9845447340485508, 0.9874977351903504, 0.960832118445256, 0.9971206498090907),
    @(1.008364823854495, 0.5401362410033301, 0.5401362410033301, 1.237696528197049, 1.063467323144168, 1.063467323144168, 1.063467323144168, 1.169027282218172, 0.9169181696287952, 1.036152779586947, 0.7345251761872527, 1.063467323144168, 1.237696528197049, 0.8889163846001896, 1.077307348912922, 0.9471000307815157, 0.8982503129430581, 0.9471000307815157, 0.9395545654933356, 0.9643371170235021, 0.9632860404775261)
)

$nRows = $rows.Length
$nCols = $rows[0].Length
$arr = New-Object 'object[,]' $nRows, $nCols
for ($i = 0; $i -lt $nRows; $i++) {
    for ($j = 0; $j -lt $nCols; $j++) {
        $arr[$i, $j] = $rows[$i][$j]
    }
}

$startRow = 4
$rng = $ws.Range($ws.Cells.Item($startRow, 3), $ws.Cells.Item($startRow + $nRows - 1, 3 + $nCols - 1))
$rng.Value2 = $arr
